# Scheduled market-data refresh for Bahamut_Profits workbook.
# Updates cached FFXIV Market Board price/profit figures (columns H-N)
# for affected Leve rows across the ALC, BSM, CRP, CUL, GSM, LTW and WVR sheets.
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#          K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 44453.652
$ws.Range("I62").Value = 60255.277
$ws.Range("J62").Value = 8900
$ws.Range("K62").Value = 60255.277
$ws.Range("L62").Value = 8900
$ws.Range("M62").Value = -59631.277
$ws.Range("N62").Value = -10148

# Row 65
$ws.Range("H65").Value = 44453.652
$ws.Range("I65").Value = 60255.277
$ws.Range("J65").Value = 8900
$ws.Range("K65").Value = 301276.385
$ws.Range("L65").Value = 44500
$ws.Range("M65").Value = -298156.385
$ws.Range("N65").Value = -50740

# Row 98
$ws.Range("H98").Value = 2084.56
$ws.Range("I98").Value = 2086.3809
$ws.Range("J98").Value = 2075
$ws.Range("K98").Value = 2086.3809
$ws.Range("L98").Value = 2075
$ws.Range("M98").Value = -588.3809000000001
$ws.Range("N98").Value = -5071

# Row 122
$ws.Range("H122").Value = 2084.56
$ws.Range("I122").Value = 2086.3809
$ws.Range("J122").Value = 2075
$ws.Range("K122").Value = 6259.1427
$ws.Range("L122").Value = 6225
$ws.Range("M122").Value = -3809.1427
$ws.Range("N122").Value = -11125

$ws = $wb.Worksheets.Item("BSM")
# Row 35
$ws.Range("H35").Value = 34750
$ws.Range("J35").Value = 34750
$ws.Range("L35").Value = 34750
$ws.Range("N35").Value = -35370

# Row 134
$ws.Range("H134").Value = 14831.707
$ws.Range("I134").Value = 1305.65
$ws.Range("J134").Value = 68935.92999999999
$ws.Range("K134").Value = 3916.95
$ws.Range("L134").Value = 206807.79
$ws.Range("M134").Value = -1381.95
$ws.Range("N134").Value = -211877.79

$ws = $wb.Worksheets.Item("CRP")
# Row 41
$ws.Range("H41").Value = 5000
$ws.Range("I41").Value = 5000
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 5000
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -4572
$ws.Range("N41").ClearContents()

# Row 99
$ws.Range("H99").Value = 1991.6364
$ws.Range("I99").Value = 1974.3158
$ws.Range("J99").Value = 2015.1428
$ws.Range("K99").Value = 1974.3158
$ws.Range("L99").Value = 2015.1428
$ws.Range("M99").Value = -476.3158000000001
$ws.Range("N99").Value = -5011.1428

# Row 126
$ws.Range("H126").Value = 1991.6364
$ws.Range("I126").Value = 1974.3158
$ws.Range("J126").Value = 2015.1428
$ws.Range("K126").Value = 5922.9474
$ws.Range("L126").Value = 6045.428400000001
$ws.Range("M126").Value = -3452.9474
$ws.Range("N126").Value = -10985.4284

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 2176.8647
$ws.Range("I5").Value = 962.5
$ws.Range("J5").Value = 2324.0605
$ws.Range("K5").Value = 2887.5
$ws.Range("L5").Value = 6972.181500000001
$ws.Range("M5").Value = -2775.5
$ws.Range("N5").Value = -7196.181500000001

# Row 39
$ws.Range("H39").Value = 19233768
$ws.Range("J39").Value = 19233768
$ws.Range("L39").Value = 57701304
$ws.Range("N39").Value = -57701892

# Row 55
$ws.Range("H55").Value = 334000
$ws.Range("I55").Value = 500500
$ws.Range("J55").Value = 1000
$ws.Range("K55").Value = 1501500
$ws.Range("L55").Value = 3000
$ws.Range("M55").Value = -1501323
$ws.Range("N55").Value = -3354

# Row 104
$ws.Range("H104").Value = 9900
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 9900
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 29700
$ws.Range("N104").Value = -34942
$ws.Range("M104").ClearContents()

# Row 122
$ws.Range("H122").Value = 27653.71
$ws.Range("I122").Value = 557.1667
$ws.Range("J122").Value = 32734.312
$ws.Range("K122").Value = 5014.5003
$ws.Range("L122").Value = 294608.808
$ws.Range("M122").Value = -2564.5003
$ws.Range("N122").Value = -299508.808

# Row 131
$ws.Range("H131").Value = 75624.96000000001
$ws.Range("I131").Value = 92260.91
$ws.Range("J131").Value = 64187.75
$ws.Range("K131").Value = 276782.73
$ws.Range("L131").Value = 192563.25
$ws.Range("M131").Value = -271742.73
$ws.Range("N131").Value = -202643.25

# Row 135
$ws.Range("H135").Value = 2176.8647
$ws.Range("I135").Value = 962.5
$ws.Range("J135").Value = 2324.0605
$ws.Range("K135").Value = 8662.5
$ws.Range("L135").Value = 20916.5445
$ws.Range("M135").Value = -6127.5
$ws.Range("N135").Value = -25986.5445

# Row 137
$ws.Range("H137").Value = 103294
$ws.Range("I137").Value = 3430
$ws.Range("J137").Value = 502750
$ws.Range("K137").Value = 10290
$ws.Range("L137").Value = 1508250
$ws.Range("M137").Value = -5190
$ws.Range("N137").Value = -1518450

$ws = $wb.Worksheets.Item("GSM")
# Row 46
$ws.Range("H46").Value = 14857.143
$ws.Range("J46").Value = 14857.143
$ws.Range("L46").Value = 14857.143
$ws.Range("N46").Value = -15169.143

# Row 57
$ws.Range("H57").Value = 23000
$ws.Range("J57").Value = 23000
$ws.Range("L57").Value = 23000
$ws.Range("N57").Value = -24640

# Row 80
$ws.Range("H80").Value = 4048.158
$ws.Range("J80").Value = 3660
$ws.Range("L80").Value = 3660
$ws.Range("N80").Value = -5656

# Row 83
$ws.Range("H83").Value = 4048.158
$ws.Range("J83").Value = 3660
$ws.Range("L83").Value = 18300
$ws.Range("N83").Value = -28284

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2673.0476
$ws.Range("I7").Value = 2544.6155
$ws.Range("J7").Value = 2881.75
$ws.Range("K7").Value = 2544.6155
$ws.Range("L7").Value = 2881.75
$ws.Range("M7").Value = -2432.6155
$ws.Range("N7").Value = -3105.75

# Row 13
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()

# Row 17
$ws.Range("H17").Value = 6673333.5
$ws.Range("I17").Value = 20000000
$ws.Range("J17").Value = 10000
$ws.Range("K17").Value = 20000000
$ws.Range("L17").Value = 10000
$ws.Range("M17").Value = -19999830
$ws.Range("N17").Value = -10340

# Row 26
$ws.Range("H26").Value = 905
$ws.Range("I26").Value = 905
$ws.Range("K26").Value = 905
$ws.Range("M26").Value = -610

# Row 30
$ws.Range("H30").Value = 7870
$ws.Range("I30").Value = 7870
$ws.Range("K30").Value = 7870
$ws.Range("M30").Value = -7762

# Row 31
$ws.Range("H31").Value = 914.3333
$ws.Range("I31").Value = 249.8
$ws.Range("J31").Value = 1745
$ws.Range("K31").Value = 249.8
$ws.Range("L31").Value = 1745
$ws.Range("M31").Value = -1.800000000000011
$ws.Range("N31").Value = -2241

# Row 34
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()

# Row 122
$ws.Range("H122").Value = 2773.0278
$ws.Range("I122").Value = 2561.0908
$ws.Range("J122").Value = 3106.0715
$ws.Range("K122").Value = 7683.2724
$ws.Range("L122").Value = 9318.2145
$ws.Range("M122").Value = -5233.2724
$ws.Range("N122").Value = -14218.2145

# Row 126
$ws.Range("H126").Value = 2673.0476
$ws.Range("I126").Value = 2544.6155
$ws.Range("J126").Value = 2881.75
$ws.Range("K126").Value = 7633.8465
$ws.Range("L126").Value = 8645.25
$ws.Range("M126").Value = -5163.8465
$ws.Range("N126").Value = -13585.25

# Row 133
$ws.Range("H133").Value = 65000
$ws.Range("J133").Value = 65000
$ws.Range("L133").Value = 65000
$ws.Range("N133").Value = -70060

# Row 136
$ws.Range("H136").Value = 4118.4224
$ws.Range("I136").Value = 2340.8
$ws.Range("J136").Value = 7673.6665
$ws.Range("K136").Value = 7022.400000000001
$ws.Range("L136").Value = 23020.9995
$ws.Range("M136").Value = -4472.400000000001
$ws.Range("N136").Value = -28120.9995

$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 2253.1765
$ws.Range("I96").Value = 2171.4285
$ws.Range("J96").Value = 2310.4
$ws.Range("K96").Value = 2171.4285
$ws.Range("L96").Value = 2310.4
$ws.Range("M96").Value = -798.4285
$ws.Range("N96").Value = -5056.4

# Row 122
$ws.Range("H122").Value = 1421.4445
$ws.Range("I122").Value = 1474.1875
$ws.Range("K122").Value = 4422.5625
$ws.Range("M122").Value = -1972.5625

# Row 126
$ws.Range("H126").Value = 1009.36365
$ws.Range("I126").Value = 990.3
$ws.Range("J126").Value = 1200
$ws.Range("K126").Value = 2970.9
$ws.Range("L126").Value = 3600
$ws.Range("M126").Value = -500.8999999999996
$ws.Range("N126").Value = -8540

# Row 132
$ws.Range("H132").Value = 4615.8667
$ws.Range("I132").Value = 1171.3334
$ws.Range("J132").Value = 18394
$ws.Range("K132").Value = 3514.0002
$ws.Range("L132").Value = 55182
$ws.Range("M132").Value = -984.0001999999999
$ws.Range("N132").Value = -60242

# Row 136
$ws.Range("H136").Value = 5178.1304
$ws.Range("I136").Value = 794.9
$ws.Range("J136").Value = 34399.668
$ws.Range("K136").Value = 2384.7
$ws.Range("L136").Value = 103199.004
$ws.Range("M136").Value = 165.3000000000002
$ws.Range("N136").Value = -108299.004
